$d = $word.ActiveDocument
$sec = $d.Sections(1)
$hdrs = $sec.Headers
$h = $hdrs.Item(1)
$h.Range.Find.Execute("Weg, Säulen, Bild, Freude", $true, $false, $false, $false, $false, $true, 1, $false, "Fluss, Weg, Regenbogen, Bild, Hand, Säulen", 2)
$h.Range.InsertParagraphAfter()
